$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-05-05 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-05-06 Monday", 2)
$d.Content.Find.Execute("285×7=", $true, $false, $false, $false, $false, $true, 1, $false, "701×7=", 2)
$d.Content.Find.Execute("181×8=", $true, $false, $false, $false, $false, $true, 1, $false, "826×6=", 2)
$d.Content.Find.Execute("484×5=", $true, $false, $false, $false, $false, $true, 1, $false, "528×7=", 2)
$d.Content.Find.Execute("998×8=", $true, $false, $false, $false, $false, $true, 1, $false, "909×2=", 2)
$d.Content.Find.Execute("162×9=", $true, $false, $false, $false, $false, $true, 1, $false, "710×4=", 2)
$d.Content.Find.Execute("430×3=", $true, $false, $false, $false, $false, $true, 1, $false, "195×6=", 2)
$d.Content.Find.Execute("920×2=", $true, $false, $false, $false, $false, $true, 1, $false, "311×5=", 2)
$d.Content.Find.Execute("781×9=", $true, $false, $false, $false, $false, $true, 1, $false, "140×7=", 2)
$d.Content.Find.Execute("222×8=", $true, $false, $false, $false, $false, $true, 1, $false, "225×4=", 2)
$d.Content.Find.Execute("517×2=", $true, $false, $false, $false, $false, $true, 1, $false, "955×2=", 2)
$d.Content.Find.Execute("997×5=", $true, $false, $false, $false, $false, $true, 1, $false, "155×4=", 2)
$d.Content.Find.Execute("215×2=", $true, $false, $false, $false, $false, $true, 1, $false, "626×9=", 2)
$d.Content.Find.Execute("117×5=", $true, $false, $false, $false, $false, $true, 1, $false, "326×6=", 2)
$d.Content.Find.Execute("891×8=", $true, $false, $false, $false, $false, $true, 1, $false, "199×6=", 2)
$d.Content.Find.Execute("578×8=", $true, $false, $false, $false, $false, $true, 1, $false, "236×6=", 2)
$d.Content.Find.Execute("824×3=", $true, $false, $false, $false, $false, $true, 1, $false, "423×4=", 2)
$d.Content.Find.Execute("913×8=", $true, $false, $false, $false, $false, $true, 1, $false, "345×2=", 2)
$d.Content.Find.Execute("109×7=", $true, $false, $false, $false, $false, $true, 1, $false, "490×4=", 2)
$d.Content.Find.Execute("845×3=", $true, $false, $false, $false, $false, $true, 1, $false, "160×5=", 2)
$d.Content.Find.Execute("239×5=", $true, $false, $false, $false, $false, $true, 1, $false, "349×6=", 2)
$d.Content.Find.Execute("995×2=", $true, $false, $false, $false, $false, $true, 1, $false, "239×9=", 2)
$d.Content.Find.Execute("123×9=", $true, $false, $false, $false, $false, $true, 1, $false, "509×8=", 2)
$d.Content.Find.Execute("426×5=", $true, $false, $false, $false, $false, $true, 1, $false, "691×6=", 2)
$d.Content.Find.Execute("800×6=", $true, $false, $false, $false, $false, $true, 1, $false, "465×8=", 2)
$d.Content.Find.Execute("400×2=", $true, $false, $false, $false, $false, $true, 1, $false, "419×5=", 2)
